# Update plots for each sample:
# Lowering/adjusting detected peak heights causes four previously
# "undetected" wildtype peaks (one per marker: CYP2D6_001, CYP2D6_003,
# CYP2D6_011, CYP2D6_013) to now be detected, which ripples into the
# allele table, the marker (genotype/phenotype) table and the final
# sample genotype result.

$wb = $excel.ActiveWorkbook

$peak    = $wb.Worksheets.Item("peak_table")
$allele  = $wb.Worksheets.Item("allele_table")
$marker  = $wb.Worksheets.Item("marker_table")
$result  = $wb.Worksheets.Item("genotype_result")

# --- peak_table: w_height (col N) drops for 4 markers -----------------
$peak.Cells.Item(2, 14).Value = 700   # CYP2D6_14 / CYP2D6_001 (S1)
$peak.Cells.Item(4, 14).Value = 600   # CYP2D6_49 / CYP2D6_003 (S1)
$peak.Cells.Item(12, 14).Value = 600  # CYP2D6_4  / CYP2D6_011 (S2)
$peak.Cells.Item(14, 14).Value = 500  # CYP2D6_17 / CYP2D6_013 (S2)

# --- allele_table: the matching wildtype allele rows now resolve a peak
# columns: K=min_height(11) M=is_detected(13) N=peak(14) O=size(15)
#          P=height(16) Q=status(17) R=message(18)

# Row 2: CYP2D6_001 / CYP2D6_14 / base G / wildtype (S1)
$allele.Cells.Item(2, 11).Value = 700
$allele.Cells.Item(2, 13).Value = $True
$allele.Cells.Item(2, 14).Value = 26
$allele.Cells.Item(2, 15).Value = 27.93
$allele.Cells.Item(2, 16).Value = 748
$allele.Cells.Item(2, 17).Value = "ok"
$allele.Cells.Item(2, 18).Value = ""

# Row 6: CYP2D6_003 / CYP2D6_49 / base T / wildtype (S1)
$allele.Cells.Item(6, 11).Value = 600
$allele.Cells.Item(6, 13).Value = $True
$allele.Cells.Item(6, 14).Value = 10
$allele.Cells.Item(6, 15).Value = 38.09
$allele.Cells.Item(6, 16).Value = 694
$allele.Cells.Item(6, 17).Value = "ok"
$allele.Cells.Item(6, 18).Value = ""

# Row 22: CYP2D6_011 / CYP2D6_4 / base G / wildtype (S2)
$allele.Cells.Item(22, 11).Value = 600
$allele.Cells.Item(22, 13).Value = $True
$allele.Cells.Item(22, 14).Value = 30
$allele.Cells.Item(22, 15).Value = 29.36
$allele.Cells.Item(22, 16).Value = 709
$allele.Cells.Item(22, 17).Value = "ok"
$allele.Cells.Item(22, 18).Value = ""

# Row 26: CYP2D6_013 / CYP2D6_17 / base C / wildtype (S2)
$allele.Cells.Item(26, 11).Value = 500
$allele.Cells.Item(26, 13).Value = $True
$allele.Cells.Item(26, 14).Value = 18
$allele.Cells.Item(26, 15).Value = 37.99
$allele.Cells.Item(26, 16).Value = 524
$allele.Cells.Item(26, 17).Value = "ok"
$allele.Cells.Item(26, 18).Value = ""

# --- marker_table: genotype (col G) / phenotype (col H) now resolved --
$marker.Cells.Item(2, 7).Value = "GG"          # CYP2D6_001
$marker.Cells.Item(2, 8).Value = "wildtype"
$marker.Cells.Item(4, 7).Value = "TT"          # CYP2D6_003
$marker.Cells.Item(4, 8).Value = "wildtype"
$marker.Cells.Item(12, 7).Value = "GG"         # CYP2D6_011
$marker.Cells.Item(12, 8).Value = "wildtype"
$marker.Cells.Item(14, 7).Value = "CC"         # CYP2D6_013
$marker.Cells.Item(14, 8).Value = "wildtype"

# --- genotype_result: final combined genotype call for the sample -----
$result.Cells.Item(2, 2).Value = "*1/*1"

# --- cosmetic leftovers from the report refresh ------------------------
$peak.Range("N14").Select()
